$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from an existing header cell (H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-15
$data = @(
    @(4, 4),
    @(4, 5),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(7, 8),
    @(7, 8),
    @(5, 6),
    @(6, 7),
    @(7, 8),
    @(5, 5),
    @(5, 6),
    @(2, 4),
    @(9, 9)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
